$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert two new columns before column D -----------------------
# This shifts the existing quarterly columns D:K -> F:M so that two new
# (more-recent) quarters can be populated in D:E.
$ws.Range("D:E").Insert()

# --- Step 2: copy number-formatting from the (now-shifted) F:G columns into
# the freshly inserted D:E columns so the new cells pick up the same date /
# number styles as the rest of the table.
$ws.Range("F:G").Copy()
$ws.Range("D:E").PasteSpecial(-4122)

# --- Step 3: populate the new D (most-recent quarter) and E (prior quarter)
# values for every data row. $null entries correspond to rows that are
# blank in the source table (section separators) and are left untouched.
$data = @{
    7 = @(43465, 43373)
    8 = @(111300, 74400)
    9 = @(68400, 43900)
    10 = @(42900, 30500)
    11 = @($null, $null)
    12 = @(1600, 1400)
    13 = @(0, 0)
    14 = @(0, 0)
    15 = @(0, 0)
    16 = @($null, $null)
    17 = @(85900, 57900)
    18 = @(25400, 16500)
    19 = @($null, $null)
    20 = @(-7300, -300)
    21 = @(19200, 17000)
    22 = @(1000, 1100)
    23 = @(17200, 15100)
    24 = @(1000, 400)
    25 = @(0, 0)
    26 = @(16200, 14700)
    27 = @(16200, 14700)
    28 = @(0, 0)
    29 = @(0, "NA")
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(7300, 300)
    33 = @(16200, 14700)
    34 = @(0, 0)
    35 = @(16200, 14700)
    38 = @(43465, 43373)
    39 = @($null, $null)
    40 = @($null, $null)
    41 = @(7100, 6200)
    42 = @(0, 0)
    43 = @(52800, 29000)
    44 = @(49500, 73300)
    45 = @(4500, 5200)
    46 = @(113800, 113700)
    47 = @(0, 0)
    48 = @(5900, 4000)
    49 = @(1000, 1100)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(1200, 700)
    53 = @(0, 0)
    54 = @(121900, 119500)
    55 = @($null, $null)
    56 = @($null, $null)
    57 = @(17700, 49000)
    58 = @(37400, 4300)
    59 = @(18500, 14600)
    60 = @(73600, 67900)
    61 = @(0, 25900)
    62 = @(10800, 2300)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(84400, 96100)
    67 = @($null, $null)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(-131500, -147700)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(37500, 23400)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(16200, 14700)
    82 = @($null, $null)
    83 = @(1100, 800)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(-1100, -3000)
    90 = @($null, $null)
    91 = @(-3000, -1600)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(-3000, -1600)
    95 = @($null, $null)
    96 = @(0, 0)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(5100, 1700)
    101 = @(-100, 0)
    102 = @(900, -2900)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($null -ne $dVal) {
        $ws.Cells.Item([int]$row, 4).Value = $dVal
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item([int]$row, 5).Value = $eVal
    }
}
